# Update beauty, MG, PSS dashboard
# The dashboard's "latest case id" cell (Sheet1!A2) moves from the
# previous run's case id to the newest one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = "CA-CCPDYWRO"
